$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old row 4 (Funcom SE drops out of the database; the former
#     row 4 AFC Ajax entry consolidates into row 3). Rows shift up automatically
#     and the sheet dimension becomes A1:AQ3. ---
$ws.Rows("4:4").Delete()

# --- Row 2: refreshed capital-structure figures; company id text stays "1" ---
$ws.Range("B2").Value = "'1"
$ws.Range("D2").Value = 0.0902
$ws.Range("E2").Value = -0.0126
$ws.Range("G2").Value = 0.04852135815991238
$ws.Range("H2").Value = 0.04852135815991238
$ws.Range("I2").Value = -0.3428258488499453
$ws.Range("J2").Value = -0.2594388354498101
$ws.Range("K2").Value = 23.3
$ws.Range("L2").Value = 0.1276013143483023
$ws.Range("M2").Value = 1.39
$ws.Range("N2").Value = 0.003849349210744946
$ws.Range("O2").Value = 0.05965665236051502
$ws.Range("P2").Value = 1.39
$ws.Range("Q2").Value = 0.003849349210744946
$ws.Range("R2").Value = 0.05965665236051502
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 25.3
$ws.Range("V2").Value = 0.07006369426751592
$ws.Range("W2").Value = 0.09769392033542977
$ws.Range("X2").Value = 0.07572721066544608
$ws.Range("Y2").Value = 0.02196670966998369
$ws.Range("Z2").Value = 1.086626637229757
$ws.Range("AA2").Value = -0.2819131493316313
$ws.Range("AB2").Value = 0.05613625286711757
$ws.Range("AC2").Value = -0.3380494021987489
$ws.Range("AD2").Value = 169.5
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 169.5
$ws.Range("AG2").Value = 144.2
$ws.Range("AH2").Value = 0.3194496796079909
$ws.Range("AI2").Value = 0.3969555035128806
$ws.Range("AJ2").Value = 0.2853750247377795
$ws.Range("AK2").Value = 0.358974358974359
$ws.Range("AL2").Value = 2.37
$ws.Range("AM2").Value = 2.098
$ws.Range("AN2").Value = 19.13092550790068
$ws.Range("AO2").Value = -26.41350210970464
$ws.Range("AP2").Value = 16.27539503386004
$ws.Range("AQ2").Value = -29.83794089609152

# --- Row 3: company renamed to AFC Ajax NV (ENXTAM:AJAX) with refreshed figures ---
$ws.Range("B3").Value = "AFC Ajax NV (ENXTAM:AJAX)"
$ws.Range("D3").Value = 0.0902
$ws.Range("E3").Value = -0.0126
$ws.Range("G3").Value = 0.04852135815991238
$ws.Range("H3").Value = 0.04852135815991238
$ws.Range("I3").Value = -0.3428258488499453
$ws.Range("J3").Value = -0.2594388354498101
$ws.Range("K3").Value = 23.3
$ws.Range("L3").Value = 0.1276013143483023
$ws.Range("M3").Value = 1.39
$ws.Range("N3").Value = 0.003849349210744946
$ws.Range("O3").Value = 0.05965665236051502
$ws.Range("P3").Value = 1.39
$ws.Range("Q3").Value = 0.003849349210744946
$ws.Range("R3").Value = 0.05965665236051502
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 25.3
$ws.Range("V3").Value = 0.07006369426751592
$ws.Range("W3").Value = 0.09769392033542977
$ws.Range("X3").Value = 0.07572721066544608
$ws.Range("Y3").Value = 0.02196670966998369
$ws.Range("Z3").Value = 1.086626637229757
$ws.Range("AA3").Value = -0.2819131493316313
$ws.Range("AB3").Value = 0.05613625286711757
$ws.Range("AC3").Value = -0.3380494021987489
$ws.Range("AD3").Value = 169.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 169.5
$ws.Range("AG3").Value = 144.2
$ws.Range("AH3").Value = 0.3194496796079909
$ws.Range("AI3").Value = 0.3969555035128806
$ws.Range("AJ3").Value = 0.2853750247377795
$ws.Range("AK3").Value = 0.358974358974359
$ws.Range("AL3").Value = 2.37
$ws.Range("AM3").Value = 2.098
$ws.Range("AN3").Value = 19.13092550790068
$ws.Range("AO3").Value = -26.41350210970464
$ws.Range("AP3").Value = 16.27539503386004
$ws.Range("AQ3").Value = -29.83794089609152

